# Fruta / hortaliza, semanal
# Insert a new weekly record at row 229 (pushing existing rows 229-247 down
# to 230-248) on the "Hortaliza, Terminal La Palmera de La Serena -
# Zapallo italiano" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 229; this shifts rows 229:247
# down to 230:248 and keeps their existing formatting/values intact.
$ws.Rows("229:229").Insert()

# Populate the newly inserted row 229 with the new weekly data point.
$ws.Range("A229").Value = 8
$ws.Range("B229").Value = "Terminal La Palmera de La Serena"
$ws.Range("C229").Value = "Coquimbo"
$ws.Range("D229").Value = 44578
$ws.Range("D229").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E229").Value = 4
$ws.Range("F229").Value = 100112032
$ws.Range("G229").Value = "Zapallo italiano"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 600
$ws.Range("K229").Value = 10500
$ws.Range("L229").Value = 11000
$ws.Range("M229").Value = 10750
$ws.Range("N229").Value = "$/caja 70 unidades"
$ws.Range("O229").Value = "Provincia de Limarí"
$ws.Range("P229").Value = 154
$ws.Range("Q229").Value = 70
$ws.Range("R229").Value = "Hortaliza"
